$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '23.124.89'
$ws.Range("E2").Value = '  -3.61%  '

$ws.Range("D3").Value = '1.601.39'
$ws.Range("E3").Value = '  -3.08%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.27'
$ws.Range("E6").Value = '  -2.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3782'
$ws.Range("E7").Value = '  -3.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3643'
$ws.Range("E8").Value = '  -4.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.92'
$ws.Range("E9").Value = '  -4.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.257'
$ws.Range("E10").Value = '  -6.91%  '

$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08140'
$ws.Range("E12").Value = '  -3.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.02'
$ws.Range("E13").Value = '  -3.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.621'
$ws.Range("E14").Value = '  -6.34%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.394'
$ws.Range("E15").Value = '  -8.01%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001254'
$ws.Range("E16").Value = '  -4.26%  '

$ws.Range("D17").Value = '1.602.31'
$ws.Range("E17").Value = '  -3.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.69'
$ws.Range("E18").Value = '  -3.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06832'
$ws.Range("E19").Value = '  -2.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.25'
$ws.Range("E20").Value = '  -7.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.567'
$ws.Range("E21").Value = '  -5.91%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.00'
$ws.Range("E23").Value = '  -5.89%  '

$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.135.71'
$ws.Range("E24").Value = '  -3.52%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.339'
$ws.Range("E25").Value = '  -4.24%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.752'
$ws.Range("E26").Value = '  -6.90%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.10'
$ws.Range("E27").Value = '  -4.52%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.52'
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.283'
$ws.Range("E29").Value = '  -2.51%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.78'
$ws.Range("E30").Value = '  -4.56%  '

$ws.Range("B31").Value = 'WEMIXTOKEN'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.426'
$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.843'
$ws.Range("E32").Value = '  -14.34%  '

$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.776.83'
$ws.Range("E33").Value = '  -2.99%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07699'
$ws.Range("E34").Value = '  -4.62%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9440'
$ws.Range("E35").Value = '  -7.45%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02746'
$ws.Range("E36").Value = '  -6.14%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.252'
$ws.Range("E37").Value = '  -7.22%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2541'
$ws.Range("E38").Value = '  -5.21%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08907'
$ws.Range("E39").Value = '  -2.27%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.08'
$ws.Range("E40").Value = '  -5.65%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.391'
$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.79'
$ws.Range("E42").Value = '  -4.63%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7108'
$ws.Range("E43").Value = '  -6.43%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.56'
$ws.Range("E44").Value = '  -4.53%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6619'
$ws.Range("E45").Value = '  -4.92%  '

$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.21%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.298'
$ws.Range("E47").Value = '  -6.49%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.974'
$ws.Range("E48").Value = '  -3.19%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.41'
$ws.Range("E49").Value = '  -1.89%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07945'
$ws.Range("E50").Value = '  -4.64%  '

$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.211'
$ws.Range("E51").Value = '  -0.98%  '
